$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-18"

# Update the header label cell (shared string "2022 (through 03-17)" -> "2022 (through 03-18)")
$ws.Range("I1").Value = "2022 (through 03-18)"

# Update March total for 2022 column (I4) and the grand Total row (I14)
$ws.Range("I4").Value = 79
$ws.Range("I14").Value = 379
